$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for the "language" field before the old "test paper" row,
# and give it the same B-column style as its neighbours.
$ws.Rows("10").Insert()
$ws.Rows("10").RowHeight = 19
$ws.Range("A10").Value = "language"
$ws.Range("B10").Value = "English"

# New annotation column (C) marking which rows are used for the online vs.
# paper test flows.
$ws.Range("C8").Value = "; online test"
$ws.Range("C9").Value = "; online test"
$ws.Range("C10").Value = "; paper test"
$ws.Range("C11").Value = "; paper test"
$ws.Range("C12").Value = "; paper test"

# The file-name cells lost their extension suffixes (the program now adds
# the extension itself), and gained an explanatory note in column D.
$ws.Range("B11").Value = "testpaper"
$ws.Range("D11").Value = "the file type (pdf) will be added by the program"

$ws.Range("B12").Value = "marksheet"
$ws.Range("D12").Value = "the file type (xlsx) will be added by the program"
